$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet: Property1 -> DataNode
$ws.Name = "DataNode"

# Update the active selection to D40 (matches the saved selection state in the target file)
[void]$ws.Range("D40").Select()
